# TS 3 Kramam Tamil Final - 19/08/2022
#
# 1) Fill in the "Observed till ????" date placeholder with "31st August 2022"
#    (with "st" as a superscript), and drop the red highlight that marked the
#    placeholder.
# 2) Reorder / tidy up the '"zlÉ" replaced with  "zgÉ" wherever applicable'
#    sentence: merge the opening quote + zlÉ + closing quote into a single
#    run, and remove the now-stray spellStart/spellEnd proofErr markers.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "????" -> "31" + superscript "st" + " August 2022"
# ---------------------------------------------------------------------------

$placeholder = $d.Content
$placeholder.Find.Execute("????", $true, $false, $false, $false, $false, $true, 1, $false, "31", 2) | Out-Null

$dateRun = $d.Content
$dateRun.Find.Execute("31", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
# The run inherited the placeholder's red highlight formatting - clear it.
$dateRun.HighlightColorIndex = 0

# Insert the superscript "st" right after "31".
$stRun = $d.Range($dateRun.End, $dateRun.End)
$stRun.InsertAfter("st")
$stRun.Font.Bold = $true
$stRun.Font.BoldBi = $true
$stRun.Font.Size = 16
$stRun.Font.SizeBi = 16
$stRun.Font.Underline = 1
$stRun.Font.Superscript = $true

# Insert " August 2022" right after "st", matching the "31" run formatting.
$restRun = $d.Range($stRun.End, $stRun.End)
$restRun.InsertAfter(" August 2022")
$restRun.Font.Bold = $true
$restRun.Font.BoldBi = $true
$restRun.Font.Size = 16
$restRun.Font.SizeBi = 16
$restRun.Font.Underline = 1

# ---------------------------------------------------------------------------
# Edit 2: rework the '"zlÉ" ... "zgÉ" ...' paragraph
# ---------------------------------------------------------------------------

$quoteParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like '*wherever applicable*') {
        $quoteParaIndex = $i
        break
    }
}

$quotePara = $d.Paragraphs($quoteParaIndex)
$quoteStart = $quotePara.Range.Start

$tailRange = $d.Range($quoteStart, $quotePara.Range.End)
$tailRange.Find.Execute("wherever applicable", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$replaceEnd = $tailRange.Start

# Clear the old run sequence (the two quoted words plus the proofErr markers
# around them), then splice in the replacement runs via WordOpenXML so the
# exact run/formatting/proofErr layout from the corrected document is
# reproduced.
$oldRange = $d.Range($quoteStart, $replaceEnd)
$oldRange.Text = ""

$newRunsXml = '<w:r><w:rPr><w:rFonts w:ascii="BRH Devanagari" w:hAnsi="BRH Devanagari" w:cs="BRH Devanagari"/><w:color w:val="000000"/><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-IN" w:eastAsia="en-IN"/></w:rPr><w:t xml:space="preserve">"zlÉ" </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/></w:rPr><w:t xml:space="preserve">replaced </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/></w:rPr><w:t>with</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BRH Malayalam" w:hAnsi="BRH Malayalam" w:cs="BRH Malayalam"/><w:color w:val="000000"/><w:sz w:val="48"/><w:szCs w:val="40"/><w:lang w:val="en-IN" w:eastAsia="en-IN"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BRH Devanagari" w:hAnsi="BRH Devanagari" w:cs="BRH Devanagari"/><w:color w:val="000000"/><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-IN" w:eastAsia="en-IN"/></w:rPr><w:t>"</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="003B3C34"><w:rPr><w:rFonts w:ascii="BRH Devanagari" w:hAnsi="BRH Devanagari" w:cs="BRH Devanagari"/><w:color w:val="000000"/><w:sz w:val="40"/><w:szCs w:val="40"/><w:highlight w:val="green"/><w:lang w:val="en-IN" w:eastAsia="en-IN"/></w:rPr><w:t>zgÉ</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BRH Devanagari" w:hAnsi="BRH Devanagari" w:cs="BRH Devanagari"/><w:color w:val="000000"/><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:val="en-IN" w:eastAsia="en-IN"/></w:rPr><w:t xml:space="preserve">" </w:t></w:r>'

$newRunsPkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $newRunsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint = $d.Range($quoteStart, $quoteStart)
$insertionPoint.InsertXML($newRunsPkg)

Write-Host "Edit complete."
